# Update factsheets with text edits from COMM
#
# The "No. of 990 Filers w/ Gov Grants" counts (column B on the detail
# sheets, column A on the Overall sheet) are converted from numeric cells
# to plain-text cells (so they render/round-trip exactly as typed, e.g. with
# a thousands separator on totals). A new statewide "Total" row is also
# appended to the County sheet, mirroring the Total row already present on
# the other breakdown sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = [string]$text
}

# ---------------------------------------------------------------------
# Overall sheet: A2 numeric 2473 -> text "2,473"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextCell $wsOverall 2 1 "2,473"

# ---------------------------------------------------------------------
# County sheet: convert the per-county counts (B2:B84) from numbers to
# text, then append a new statewide Total row (row 85).
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
$countyCounts = @(2,2,20,19,5,3,3,12,21,7,54,4,40,8,19,5,15,3,10,1,9,8,33,20,83,2,6,47,9,7,15,9,150,10,6,7,18,34,109,1,214,2,3,11,16,23,20,1,10,111,8,26,10,10,5,42,3,21,6,2,46,10,239,6,6,2,4,1,17,87,2,4,58,6,1,9,32,10,3,13,167,361,9)
for ($i = 0; $i -lt $countyCounts.Length; $i++) {
    $row = $i + 2
    Set-TextCell $wsCounty $row 2 $countyCounts[$i]
}

$totalRow = 85
Set-TextCell $wsCounty $totalRow 1 "Total"
Set-TextCell $wsCounty $totalRow 2 "2,473"
Set-TextCell $wsCounty $totalRow 3 '$4,577,112,023'
Set-TextCell $wsCounty $totalRow 4 "9.00%"
Set-TextCell $wsCounty $totalRow 5 "-10.23%"
Set-TextCell $wsCounty $totalRow 6 "65.79%"

# ---------------------------------------------------------------------
# Congressional District sheet: B2:B14 numbers -> text; Total row (B15)
# numeric 2473 -> text "2,473"
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
$cdCounts = @(298,115,142,144,240,134,269,272,136,209,223,203,88)
for ($i = 0; $i -lt $cdCounts.Length; $i++) {
    $row = $i + 2
    Set-TextCell $wsCd $row 2 $cdCounts[$i]
}
Set-TextCell $wsCd 15 2 "2,473"

# ---------------------------------------------------------------------
# Size sheet: B2:B7 numbers -> text; Total row (B8) numeric 2473 -> text
# "2,473"
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @(735,687,459,169,319,104)
for ($i = 0; $i -lt $sizeCounts.Length; $i++) {
    $row = $i + 2
    Set-TextCell $wsSize $row 2 $sizeCounts[$i]
}
Set-TextCell $wsSize 8 2 "2,473"

# ---------------------------------------------------------------------
# Subsector sheet: B2:B13 numbers -> text; Total row (B14) numeric 2473
# -> text "2,473"
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
$subCounts = @(209,174,117,211,49,821,20,2,202,67,575,26)
for ($i = 0; $i -lt $subCounts.Length; $i++) {
    $row = $i + 2
    Set-TextCell $wsSub $row 2 $subCounts[$i]
}
Set-TextCell $wsSub 14 2 "2,473"
